# feat: add 2022-Q3 data
#
# The existing "2022-Q2" sheet is duplicated to create a new "2022-Q3" sheet
# (placed right before the original, so tab order becomes 总计, 2022-Q3,
# 2022-Q2), its figures are updated to the Q3 numbers, and the "总计" summary
# sheet gets a new top data row for 2022-Q3 (the old 2022-Q2 row shifts down).

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "2022-Q2" into a new "2022-Q3" sheet, inserted before it ---
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2, [System.Reflection.Missing]::Value)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# --- 2. Update the new Q3 sheet's data row with the Q3 figures ---
# D2:G2 hold numeric-looking values stored as text (matching the source
# file's convention), so force text via NumberFormat before assigning, then
# drop the format again so no stray style sticks to the cells.
$wsQ3.Range("D2:G2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "15.28"
$wsQ3.Range("E2").Value = "77.14"
$wsQ3.Range("F2").Value = "0.86"
$wsQ3.Range("G2").Value = "0.1314"
$wsQ3.Range("D2:G2").ClearFormats()
$wsQ3.Range("H2").Value = 6

# --- 3. Insert a new 2022-Q3 row into the "总计" summary sheet ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# The row insert leaves stray formatting on B2:D2 picked up from the shift;
# clear it so they end up styleless, matching the original row's cells.
$wsTotal.Range("B2:D2").ClearFormats()

# A2 needs the same style as A3 (bold/centered/bordered) - copy it over via
# paste-special (formats only) rather than Range.Style, which didn't stick.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.13

$wsTotal.Range("A3").Value = 1
